# refatoracao - calculos de apoio medio
# Insert the new "apoio" detail columns (std/min/max) right after apoio_medio (col L),
# then insert the new "contribuicoes" detail columns (std/min/max) right after the
# renamed contribuicoes_med column (which lands at col Q once the first insert has
# shifted things over). Inserting (rather than just writing into blank cells) makes
# Excel copy the left neighbour's style onto the new cells/column headers for free.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) apoio_std / apoio_min / apoio_max go in right after apoio_medio (L)
$ws.Columns("M:O").Insert(-4161)

# 2) contribuicoes_std / contribuicoes_min / contribuicoes_max go in right after
#    contribuicoes_med, which (post step 1) sits at column Q.
$ws.Columns("R:T").Insert(-4161)

# ---- Header row (row 1) ----
$ws.Range("H1").Value = "arrecadado_avg"
$ws.Range("I1").Value = "arrecadado_std"
$ws.Range("J1").Value = "arrecadado_min"
$ws.Range("K1").Value = "arrecadado_max"
$ws.Range("M1").Value = "apoio_std"
$ws.Range("N1").Value = "apoio_min"
$ws.Range("O1").Value = "apoio_max"
$ws.Range("Q1").Value = "contribuicoes_med"
$ws.Range("R1").Value = "contribuicoes_std"
$ws.Range("S1").Value = "contribuicoes_min"
$ws.Range("T1").Value = "contribuicoes_max"

# ---- Data rows ----
# Row 2 (coletivo)
$ws.Range("L2").Value = 94.9905854649608
$ws.Range("M2").Value = 47.6922463509549
$ws.Range("N2").Value = 47.35034461927121
$ws.Range("O2").Value = 305.2480444061168
$ws.Range("R2").Value = 212.2582078460797
$ws.Range("S2").Value = 35
$ws.Range("T2").Value = 808

# Row 3 (empresa)
$ws.Range("L3").Value = 110.6538302669828
$ws.Range("M3").Value = 45.12744090621267
$ws.Range("N3").Value = 39.22956647121969
$ws.Range("O3").Value = 257.7853211115706
$ws.Range("R3").Value = 378.1809353534696
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 1711

# Row 4 (feminino)
$ws.Range("L4").Value = 82.69413375987617
$ws.Range("M4").Value = 30.90119243508478
$ws.Range("N4").Value = 13.93896149503088
$ws.Range("O4").Value = 194.2230576381307
$ws.Range("R4").Value = 547.4955526904555
$ws.Range("S4").Value = 3
$ws.Range("T4").Value = 5879

# Row 5 (masculino)
$ws.Range("L5").Value = 91.3503645951285
$ws.Range("M5").Value = 52.58131393014926
$ws.Range("N5").Value = 21.61624650544615
$ws.Range("O5").Value = 792.0360759681182
$ws.Range("R5").Value = 401.4417134786221
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = 6494

# Row 6 (outros)
$ws.Range("L6").Value = 53.14416408875834
$ws.Range("M6").Value = 7.563317519432532
$ws.Range("N6").Value = 47.79609098250058
$ws.Range("O6").Value = 58.4922371950161
$ws.Range("R6").Value = 18.38477631085023
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 36
